# Update "handback" timestamp cells as part of regenerating the report.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file.
$wsOverview.Range("G2").Value = "2016-08-19 01:02:58"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file.
$wsZhCn.Range("H2").Value = "2016-08-19 01:02:53"
$wsZhCn.Range("K2").Value = "2016-08-19 01:03:15"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file.
$wsDeDe.Range("H2").Value = "2016-08-19 01:02:58"
$wsDeDe.Range("K2").Value = "2016-08-19 01:03:22"
